$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 34.22762
$ws.Range("H2").Value = 102.68286
$ws.Range("I2").Value = 0.2984485907090856
$ws.Range("J2").Value = 0.3095564710510569
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1613523333333333
$ws.Range("N2").Value = 0.484057
$ws.Range("O2").Value = 0.2926766298022186
$ws.Range("P2").Value = 0.3782017374917083
$ws.Range("Q2").Value = 5.522706351446668
$ws.Range("R2").Value = 49.70435716302001
$ws.Range("S2").Value = 0.08734892769795694
$ws.Range("T2").Value = 0.1170747952033114
$ws.Range("G3").Value = 34.22762
$ws.Range("H3").Value = 102.68286
$ws.Range("I3").Value = 0.2984485907090856
$ws.Range("J3").Value = 0.3095564710510569
$ws.Range("O3").Value = 0.02891473894686308
$ws.Range("P3").Value = 0.03736411928828315
$ws.Range("Q3").Value = 0.5456110812133335
$ws.Range("R3").Value = 4.910499730920001
$ws.Range("S3").Value = 0.008629563089412397
$ws.Range("T3").Value = 0.01156630491081166
$ws.Range("G4").Value = 34.22762
$ws.Range("H4").Value = 102.68286
$ws.Range("I4").Value = 0.2984485907090856
$ws.Range("J4").Value = 0.3095564710510569
$ws.Range("M4").Value = 0.3740059999999999
$ws.Range("N4").Value = 0.7480119999999999
$ws.Range("O4").Value = 0.6784086312509182
$ws.Range("P4").Value = 0.5844341432200085
$ws.Range("Q4").Value = 12.80133524572
$ws.Range("R4").Value = 76.80801147432
$ws.Range("S4").Value = 0.2024700999217163
$ws.Range("T4").Value = 0.1809153709369338
$ws.Range("G5").Value = 65.67978099999999
$ws.Range("I5").Value = 0.5726964970842663
$ws.Range("J5").Value = 0.594011538803056
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1613523333333333
$ws.Range("N5").Value = 0.484057
$ws.Range("O5").Value = 0.2926766298022186
$ws.Range("P5").Value = 0.3782017374917083
$ws.Range("Q5").Value = 10.59758591717233
$ws.Range("R5").Value = 95.37827325455099
$ws.Range("S5").Value = 0.1676148806661592
$ws.Range("T5").Value = 0.2246561960654391
$ws.Range("G6").Value = 65.67978099999999
$ws.Range("I6").Value = 0.5726964970842663
$ws.Range("J6").Value = 0.594011538803056
$ws.Range("O6").Value = 0.02891473894686308
$ws.Range("P6").Value = 0.03736411928828315
$ws.Range("R6").Value = 9.422815460946
$ws.Range("S6").Value = 0.01655936970897449
$ws.Range("T6").Value = 0.02219471799445402
$ws.Range("G7").Value = 65.67978099999999
$ws.Range("I7").Value = 0.5726964970842663
$ws.Range("J7").Value = 0.594011538803056
$ws.Range("M7").Value = 0.3740059999999999
$ws.Range("N7").Value = 0.7480119999999999
$ws.Range("O7").Value = 0.6784086312509182
$ws.Range("P7").Value = 0.5844341432200085
$ws.Range("Q7").Value = 24.56463217268599
$ws.Range("R7").Value = 147.387793036116
$ws.Range("S7").Value = 0.3885222467091325
$ws.Range("T7").Value = 0.3471606247431629
$ws.Range("G8").Value = 0.871228
$ws.Range("H8").Value = 2.613684
$ws.Range("I8").Value = 0.007596694388517088
$ws.Range("J8").Value = 0.007879433777775674
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1613523333333333
$ws.Range("N8").Value = 0.484057
$ws.Range("O8").Value = 0.2926766298022186
$ws.Range("P8").Value = 0.3782017374917083
$ws.Range("Q8").Value = 0.1405746706653334
$ws.Range("R8").Value = 1.265172035988
$ws.Range("S8").Value = 0.002223374911268607
$ws.Range("T8").Value = 0.002980015545205615
$ws.Range("G9").Value = 0.871228
$ws.Range("H9").Value = 2.613684
$ws.Range("I9").Value = 0.007596694388517088
$ws.Range("J9").Value = 0.007879433777775674
$ws.Range("O9").Value = 0.02891473894686308
$ws.Range("P9").Value = 0.03736411928828315
$ws.Range("Q9").Value = 0.01388795513866667
$ws.Range("R9").Value = 0.124991596248
$ws.Range("S9").Value = 0.0002196564351030712
$ws.Range("T9").Value = 0.0002944081035969378
$ws.Range("G10").Value = 0.871228
$ws.Range("H10").Value = 2.613684
$ws.Range("I10").Value = 0.007596694388517088
$ws.Range("J10").Value = 0.007879433777775674
$ws.Range("M10").Value = 0.3740059999999999
$ws.Range("N10").Value = 0.7480119999999999
$ws.Range("O10").Value = 0.6784086312509182
$ws.Range("P10").Value = 0.5844341432200085
$ws.Range("Q10").Value = 0.325844499368
$ws.Range("R10").Value = 1.955066996208
$ws.Range("S10").Value = 0.005153663042145409
$ws.Range("T10").Value = 0.004605010128973121
$ws.Range("G11").Value = 12.345814
$ws.Range("H11").Value = 24.691628
$ws.Range("I11").Value = 0.1076496346943346
$ws.Range("J11").Value = 0.07443747893451221
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1613523333333333
$ws.Range("N11").Value = 0.484057
$ws.Range("O11").Value = 0.2926766298022186
$ws.Range("P11").Value = 0.3782017374917083
$ws.Range("Q11").Value = 1.992025895799334
$ws.Range("R11").Value = 11.952155374796
$ws.Range("S11").Value = 0.03150653228177783
$ws.Range("T11").Value = 0.02815238386753495
$ws.Range("G12").Value = 12.345814
$ws.Range("H12").Value = 24.691628
$ws.Range("I12").Value = 0.1076496346943346
$ws.Range("J12").Value = 0.07443747893451221
$ws.Range("O12").Value = 0.02891473894686308
$ws.Range("P12").Value = 0.03736411928828315
$ws.Range("Q12").Value = 0.1968005057026667
$ws.Range("R12").Value = 1.180803034216
$ws.Range("S12").Value = 0.003112661084911858
$ws.Range("T12").Value = 0.002781290842428178
$ws.Range("G13").Value = 12.345814
$ws.Range("H13").Value = 24.691628
$ws.Range("I13").Value = 0.1076496346943346
$ws.Range("J13").Value = 0.07443747893451221
$ws.Range("M13").Value = 0.3740059999999999
$ws.Range("N13").Value = 0.7480119999999999
$ws.Range("O13").Value = 0.6784086312509182
$ws.Range("P13").Value = 0.5844341432200085
$ws.Range("Q13").Value = 4.617408510883999
$ws.Range("R13").Value = 18.469634043536
$ws.Range("S13").Value = 0.07303044132764486
$ws.Range("T13").Value = 0.04350380422454907
$ws.Range("G14").Value = 1.560702333333333
$ws.Range("H14").Value = 4.682107
$ws.Range("I14").Value = 0.01360858312379636
$ws.Range("J14").Value = 0.01411507743359944
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1613523333333333
$ws.Range("N14").Value = 0.484057
$ws.Range("O14").Value = 0.2926766298022186
$ws.Range("P14").Value = 0.3782017374917083
$ws.Range("Q14").Value = 0.2518229631221112
$ws.Range("R14").Value = 2.266406668099
$ws.Range("S14").Value = 0.003982914245056069
$ws.Range("T14").Value = 0.005338346810217312
$ws.Range("G15").Value = 1.560702333333333
$ws.Range("H15").Value = 4.682107
$ws.Range("I15").Value = 0.01360858312379636
$ws.Range("J15").Value = 0.01411507743359944
$ws.Range("O15").Value = 0.02891473894686308
$ws.Range("P15").Value = 0.03736411928828315
$ws.Range("Q15").Value = 0.02487863566155556
$ws.Range("R15").Value = 0.223907720954
$ws.Range("S15").Value = 0.0003934886284612583
$ws.Range("T15").Value = 0.0005273974369923632
$ws.Range("G16").Value = 1.560702333333333
$ws.Range("H16").Value = 4.682107
$ws.Range("I16").Value = 0.01360858312379636
$ws.Range("J16").Value = 0.01411507743359944
$ws.Range("M16").Value = 0.3740059999999999
$ws.Range("N16").Value = 0.7480119999999999
$ws.Range("O16").Value = 0.6784086312509182
$ws.Range("P16").Value = 0.5844341432200085
$ws.Range("Q16").Value = 0.5837120368806666
$ws.Range("R16").Value = 3.502272221284
$ws.Range("S16").Value = 0.009232180250279037
$ws.Range("T16").Value = 0.008249333186389766
